$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update S2:S6 to the new total_egg_deposition value (sum of R2:R6, the Female fecundity totals)
$ws.Range("S2:S6").Value = 21949190.17211484

# M7 changes to match E7 (43/55)
$ws.Range("M7").Value = 0.7818181818181819

# Clear the now-unused S column cells for rows 7 through 26
$ws.Range("S7:S26").ClearContents()
